# Auto-generated edit script applying the Raiden_Profits.xlsx value refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 3454.7693
$ws.Range("I2").Value = 1924
$ws.Range("J2").Value = 4135.1113
$ws.Range("K2").Value = 1924
$ws.Range("L2").Value = 4135.1113
$ws.Range("M2").Value = -1811
$ws.Range("N2").Value = -4361.1113
# Row 125
$ws.Range("H125").Value = 3231.2354
$ws.Range("I125").Value = 10106.5
$ws.Range("J125").Value = 1115.7693
$ws.Range("K125").Value = 90958.5
$ws.Range("L125").Value = 10041.9237
$ws.Range("M125").Value = -88498.5
$ws.Range("N125").Value = -14961.9237

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 520.6667
$ws.Range("I5").Value = 533
$ws.Range("J5").Value = 496
$ws.Range("K5").Value = 533
$ws.Range("L5").Value = 496
$ws.Range("M5").Value = -421
$ws.Range("N5").Value = -720
# Row 61
$ws.Range("H61").Value = 2328.7144
$ws.Range("I61").Value = 1643.6923
$ws.Range("K61").Value = 1643.6923
$ws.Range("M61").Value = -1431.6923
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
# Row 132
$ws.Range("H132").Value = 3837.7778
$ws.Range("I132").Value = 3923.8333
$ws.Range("J132").Value = 3665.6667
$ws.Range("K132").Value = 11771.4999
$ws.Range("L132").Value = 10997.0001
$ws.Range("M132").Value = -9241.499899999999
$ws.Range("N132").Value = -16057.0001
# Row 136
$ws.Range("H136").Value = 2328.7144
$ws.Range("I136").Value = 1643.6923
$ws.Range("K136").Value = 4931.0769
$ws.Range("M136").Value = -2381.0769

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 520.6667
$ws.Range("I4").Value = 533
$ws.Range("J4").Value = 496
$ws.Range("K4").Value = 533
$ws.Range("L4").Value = 496
$ws.Range("M4").Value = -418
$ws.Range("N4").Value = -726
# Row 105
$ws.Range("H105").Value = 2193.2856
$ws.Range("I105").Value = 1635.2632
$ws.Range("K105").Value = 1635.2632
$ws.Range("M105").Value = 111.7367999999999
# Row 126
$ws.Range("H126").Value = 42890
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 42890
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 42890
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -52770
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 38
$ws.Range("H38").Value = 7025
$ws.Range("J38").Value = 10499.5
$ws.Range("L38").Value = 10499.5
$ws.Range("N38").Value = -11253.5
# Row 46
$ws.Range("H46").Value = 7025
$ws.Range("J46").Value = 10499.5
$ws.Range("L46").Value = 10499.5
$ws.Range("N46").Value = -10921.5
# Row 62
$ws.Range("H62").Value = 15336.16
$ws.Range("I62").Value = 13533.286
$ws.Range("K62").Value = 13533.286
$ws.Range("M62").Value = -12909.286
# Row 65
$ws.Range("H65").Value = 15336.16
$ws.Range("I65").Value = 13533.286
$ws.Range("K65").Value = 67666.42999999999
$ws.Range("M65").Value = -64546.42999999999
# Row 98
$ws.Range("H98").Value = 50780
$ws.Range("J98").Value = 50780
$ws.Range("L98").Value = 50780
$ws.Range("N98").Value = -55272
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
# Row 132
$ws.Range("H132").Value = 947.7083
$ws.Range("I132").Value = 954.6087
$ws.Range("K132").Value = 2863.8261
$ws.Range("M132").Value = -333.8261000000002

$ws = $wb.Worksheets.Item("CUL")
# Row 60
$ws.Range("H60").Value = 433.55554
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 98
$ws.Range("H98").Value = 289.875
$ws.Range("I98").Value = 280
$ws.Range("K98").Value = 840
$ws.Range("M98").Value = 658
# Row 113
$ws.Range("H113").Value = 992.875
$ws.Range("J113").Value = 1149.4
$ws.Range("L113").Value = 3448.2
$ws.Range("N113").Value = -7788.200000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 109.14286
$ws.Range("I2").Value = 109.14286
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 109.14286
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 3.857140000000001
$ws.Range("N2").ClearContents()
# Row 122
$ws.Range("H122").Value = 2782
$ws.Range("I122").Value = 2782
$ws.Range("K122").Value = 8346
$ws.Range("M122").Value = -5896
# Row 127
$ws.Range("H127").Value = 49998.5
$ws.Range("I127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("M127").ClearContents()
# Row 132
$ws.Range("H132").Value = 1945.2354
$ws.Range("I132").Value = 2397.3333
$ws.Range("J132").Value = 1436.625
$ws.Range("K132").Value = 7191.999899999999
$ws.Range("L132").Value = 4309.875
$ws.Range("M132").Value = -4661.999899999999
$ws.Range("N132").Value = -9369.875
# Row 135
$ws.Range("H135").Value = 67500
$ws.Range("J135").Value = 67500
$ws.Range("L135").Value = 67500
$ws.Range("N135").Value = -77640

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3636.7144
$ws.Range("I40").Value = 3684.6086
$ws.Range("K40").Value = 3684.6086
$ws.Range("M40").Value = -3548.6086
# Row 46
$ws.Range("H46").Value = 1103.8667
$ws.Range("I46").Value = 528.6667
$ws.Range("J46").Value = 1247.6666
$ws.Range("K46").Value = 528.6667
$ws.Range("L46").Value = 1247.6666
$ws.Range("M46").Value = -340.6667
$ws.Range("N46").Value = -1623.6666
# Row 132
$ws.Range("H132").Value = 2916.7144
$ws.Range("I132").Value = 2883.6
$ws.Range("K132").Value = 8650.799999999999
$ws.Range("M132").Value = -6120.799999999999
# Row 136
$ws.Range("H136").Value = 1506.9231
$ws.Range("I136").Value = 1169.1
$ws.Range("K136").Value = 3507.3
$ws.Range("M136").Value = -957.2999999999997

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 14200.25
$ws.Range("I2").Value = 14900
$ws.Range("J2").Value = 13500.5
$ws.Range("K2").Value = 14900
$ws.Range("L2").Value = 13500.5
$ws.Range("M2").Value = -14788
$ws.Range("N2").Value = -13724.5
# Row 4
$ws.Range("H4").Value = 3967.875
$ws.Range("J4").Value = 3471.2856
$ws.Range("L4").Value = 3471.2856
$ws.Range("N4").Value = -3697.2856
# Row 107
$ws.Range("H107").Value = 1367.3158
$ws.Range("I107").Value = 1404.1333
$ws.Range("J107").Value = 1343.3043
$ws.Range("K107").Value = 4212.3999
$ws.Range("L107").Value = 4029.9129
$ws.Range("M107").Value = -2292.3999
$ws.Range("N107").Value = -7869.9129
# Row 122
$ws.Range("H122").Value = 7322.381
$ws.Range("I122").Value = 7893.2104
$ws.Range("J122").Value = 1899.5
$ws.Range("K122").Value = 23679.6312
$ws.Range("L122").Value = 5698.5
$ws.Range("M122").Value = -21229.6312
$ws.Range("N122").Value = -10598.5
# Row 132
$ws.Range("H132").Value = 2511.4583
$ws.Range("I132").Value = 2577.3845
$ws.Range("J132").Value = 2225.7778
$ws.Range("K132").Value = 7732.1535
$ws.Range("L132").Value = 6677.3334
$ws.Range("M132").Value = -5202.1535
$ws.Range("N132").Value = -11737.3334
# Row 136
$ws.Range("H136").Value = 639.3125
$ws.Range("I136").Value = 444.84616
$ws.Range("J136").Value = 1482
$ws.Range("K136").Value = 1334.53848
$ws.Range("L136").Value = 4446
$ws.Range("M136").Value = 1215.46152
$ws.Range("N136").Value = -9546

Write-Host "Applied 190 cell updates across 8 sheets"